# Update crypto price/volume table values per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.119.37"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "2.549.51"
$ws.Range("E4").Value = "  -0.07%  "
$cell = $ws.Range("D5")
$cell.Value = "'539.78"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "
$cell = $ws.Range("D6")
$cell.Value = "'144.14"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("E7").Value = "  -0.12%  "
$cell = $ws.Range("D8")
$cell.Value = "'0.571"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "2.566.04"
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("E11").Value = "  +1.72%  "
$cell = $ws.Range("D12")
$cell.Value = "'5.47"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("E13").Value = "  +3.75%  "
$ws.Range("D14").Value = "2.998.13"
$ws.Range("E14").Value = "  +2.40%  "
$cell = $ws.Range("D15")
$cell.Value = "'24.08"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "60.060.83"
$ws.Range("E16").Value = "  +2.55%  "
$cell = $ws.Range("D17")
$cell.Value = "'0.0000144"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +4.53%  "
$ws.Range("D18").Value = "2.570.46"
$ws.Range("E18").Value = "  +2.70%  "
$ws.Range("E19").Value = "  -0.81%  "
$cell = $ws.Range("D20")
$cell.Value = "'4.33"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.52%  "
$cell = $ws.Range("D21")
$cell.Value = "'327.06"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("E22").Value = "  +0.37%  "
$cell = $ws.Range("D23")
$cell.Value = "'5.97"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +3.99%  "
$cell = $ws.Range("D24")
$cell.Value = "'63.29"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +4.24%  "
$ws.Range("E25").Value = "  -0.49%  "
$cell = $ws.Range("D26")
$cell.Value = "'0.166"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +3.79%  "
$cell = $ws.Range("D27")
$cell.Value = "'0.994"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.19%  "
$cell = $ws.Range("D28")
$cell.Value = "'8.05"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +4.42%  "
$cell = $ws.Range("D29")
$cell.Value = "'7.07"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +3.91%  "
$ws.Range("D30").Value = "0.0₃0795"
$ws.Range("E30").Value = "  +4.11%  "
$cell = $ws.Range("D31")
$cell.Value = "'1.82"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +2.20%  "
$ws.Range("E32").Value = "  -3.91%  "
$cell = $ws.Range("D33")
$cell.Value = "'165.02"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +4.89%  "
$ws.Range("E34").Value = "  +5.30%  "
$ws.Range("E35").Value = "  +0.25%  "
$cell = $ws.Range("D36")
$cell.Value = "'18.75"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("E38").Value = "  +2.29%  "
$cell = $ws.Range("D39")
$cell.Value = "'37.03"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +1.01%  "
$cell = $ws.Range("D40")
$cell.Value = "'303.62"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.88%  "
$cell = $ws.Range("D41")
$cell.Value = "'5.61"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -5.01%  "
$ws.Range("E42").Value = "  +6.17%  "
$ws.Range("E43").Value = "  +2.16%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$cell = $ws.Range("D44")
$cell.Value = "'0.611"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$cell = $ws.Range("D45")
$cell.Value = "'0.994"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.06%  "
$cell = $ws.Range("D46")
$cell.Value = "'10.85"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.77%  "
$cell = $ws.Range("D47")
$cell.Value = "'127.42"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +2.47%  "
$cell = $ws.Range("D48")
$cell.Value = "'0.0938"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("E51").Value = "  +1.11%  "
